$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update points (column I) values for both pursuit blocks (rows 2-31 and 32-61)
# included stuff for relays: new scoring scale
$ws.Range("I2").Value = 400
$ws.Range("I3").Value = 320
$ws.Range("I4").Value = 240
$ws.Range("I5").Value = 200
$ws.Range("I6").Value = 180
$ws.Range("I7").Value = 160
$ws.Range("I8").Value = 144
$ws.Range("I9").Value = 128
$ws.Range("I10").Value = 116
$ws.Range("I11").Value = 104
$ws.Range("I12").Value = 96
$ws.Range("I13").Value = 88
$ws.Range("I14").Value = 80
$ws.Range("I15").Value = 72
$ws.Range("I16").Value = 64
$ws.Range("I17").Value = 60
$ws.Range("I18").Value = 56
$ws.Range("I19").Value = 52
$ws.Range("I20").Value = 48
$ws.Range("I21").Value = 44
$ws.Range("I22").Value = 40
$ws.Range("I23").Value = 36
$ws.Range("I24").Value = 32
$ws.Range("I25").Value = 28
$ws.Range("I26").Value = 24
$ws.Range("I27").Value = 20
$ws.Range("I28").Value = 20
$ws.Range("I29").Value = 20
$ws.Range("I30").Value = 20
$ws.Range("I31").Value = 20
$ws.Range("I32").Value = 400
$ws.Range("I33").Value = 320
$ws.Range("I34").Value = 240
$ws.Range("I35").Value = 200
$ws.Range("I36").Value = 180
$ws.Range("I37").Value = 160
$ws.Range("I38").Value = 144
$ws.Range("I39").Value = 128
$ws.Range("I40").Value = 116
$ws.Range("I41").Value = 104
$ws.Range("I42").Value = 96
$ws.Range("I43").Value = 88
$ws.Range("I44").Value = 80
$ws.Range("I45").Value = 72
$ws.Range("I46").Value = 64
$ws.Range("I47").Value = 60
$ws.Range("I48").Value = 56
$ws.Range("I49").Value = 52
$ws.Range("I50").Value = 48
$ws.Range("I51").Value = 44
$ws.Range("I52").Value = 40
$ws.Range("I53").Value = 36
$ws.Range("I54").Value = 32
$ws.Range("I55").Value = 28
$ws.Range("I56").Value = 24
$ws.Range("I57").Value = 20
$ws.Range("I58").Value = 20
$ws.Range("I59").Value = 20
$ws.Range("I60").Value = 20
$ws.Range("I61").Value = 20
